$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (C3:H3): give it the new "themed" fill (bold font + border already present) ---
# Interior.Color set to the RGB equivalent of Theme "Blue, Accent 1, Lighter 80%"
# (the engine does not persist Interior.TintAndShade, so the plain RGB value is the closest
# achievable reproduction of fgColor theme="4" tint="0.8").
$ws.Range("C3:H3").Interior.Color = 16247774

# --- Column C: widen + drop the old best-fit flag ---
# ColumnWidth is internally offset by Excel's default padding (~5/6 of a character) before
# being written out, so we back that off here to land on the target stored width.
$ws.Columns("C").ColumnWidth = 35.5

# --- New data row (row 10) ---
$ws.Range("C10").Value = "C021_FT Dataprep_Speed Profile.ipynb"
$ws.Range("D10").Value = "/code-cloud/"
$ws.Range("E10").Value = "raw_base_2023-06-05.csv"
$ws.Range("F10").Value = "f'dataprep_speed_profile_{todaydt}.csv'"
$ws.Range("G10").Value = "Create speed profile variables such as `ndistance total, distance group, race time group etc`n"

# G10 picks up the same wrap-text/border style already used by the description column elsewhere.
$ws.Range("G10").WrapText = $true

# Row 10 grows to fit the wrapped description (matches the height Excel would compute for it).
$ws.Rows(10).RowHeight = 57.6

# --- View state: put the cursor on C10 and scroll the new row into view ---
$ws.Range("C10").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
